$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") — copy the header formatting
# (bold font, thin border, centered/top alignment) from the existing H1
# header cell so the new headers reuse the same style as the rest of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data columns I (I0) and J (IF), rows 2-11.
$i0Values = @(5, 8, 2, 8, 5, 5, 9, 7, 9, 6)
$ifValues = @(5, 8, 3, 9, 5, 5, 9, 7, 9, 6)

for ($idx = 0; $idx -lt $i0Values.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $i0Values[$idx]
    $ws.Cells.Item($row, 10).Value = $ifValues[$idx]
}
